# ---------------------------------------------------------------------------
# Insert the "2022-Q3" worksheet (fund-holding detail) right after "总计",
# and add the corresponding summary row to "总计".  Matches the commit
# "feat: add 2022-Q3 data".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Helper-ish inline pattern used throughout:
#   1) write the value with a leading "'" so Excel stores it as TEXT
#      (otherwise a numeric-looking string like "42.79" becomes a Number).
#   2) if the target cell should carry the bold/bordered "header/index"
#      style, re-apply that exact style by Copy + PasteSpecial(Formats)
#      from a cell that already has it -- doing this *after* the value
#      write overwrites the transient quote-prefix style Excel created.
#   3) if the target cell should have NO special style at all, reset it
#      with Style = "Normal" (clears the quote-prefix flag Excel added).
# ---------------------------------------------------------------------------

# ---- 1. Create the new sheet right after "总计" -----------------------
$zongji = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Style donors (cells that already carry the workbook's header / index style)
$styleDonorSheet = $wb.Worksheets.Item("2022-Q2")
$headerStyleDonor = $styleDonorSheet.Range("B1")
$indexStyleDonor = $styleDonorSheet.Range("A2")

# ---- 2. Header row (B1:H1) ------------------------------------------------
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = "'" + $headers[$i]
    $headerStyleDonor.Copy()
    $cell.PasteSpecial($xlPasteFormats)
}

# ---- 3. Data rows (rows 2..17) -------------------------------------------
$rows = @(
    @("001532","华安文体健康主题灵活配置混合A","42.79","87.37","2.75","1.1767",4),
    @("002350","华安安华灵活配置混合A","32.74","87.55","2.68","0.8774",3),
    @("519198","万家颐和灵活配置混合A","9.55","93.93","8.51","0.8127",2),
    @("014207","华安产业精选混合A","26.65","82.10","2.64","0.7036",2),
    @("014208","华安产业精选混合C","22.90","82.10","2.64","0.6046",2),
    @("008979","万家民丰回报一年持有期混合","18.51","29.35","1.94","0.3591",8),
    @("013680","华安品质甄选混合A","10.86","42.70","1.82","0.1977",2),
    @("519197","万家颐达灵活配置混合","2.25","45.36","5.98","0.1346",3),
    @("013116","华安文体健康主题灵活配置混合C","4.16","87.37","2.75","0.1144",4),
    @("013681","华安品质甄选混合C","4.45","42.70","1.82","0.0810",2),
    @("010690","万家互联互通核心资产量化策略混合A","0.55","92.41","7.10","0.0390",7),
    @("011629","银河核心优势混合","2.68","26.95","1.10","0.0295",9),
    @("001267","泰达宏利蓝筹价值混合","0.49","94.19","4.40","0.0216",9),
    @("016620","万家颐和灵活配置混合C","0.18","93.93","8.51","0.0153",2),
    @("016183","华安安华灵活配置混合C","0.45","87.55","2.68","0.0121",3),
    @("010691","万家互联互通核心资产量化策略混合C","0.15","92.41","7.10","0.0106",7)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $rowData = $rows[$r]

    # Column A: numeric row index, styled like the other quarter sheets.
    $q3.Cells.Item($rowNum, 1).Value = $r
    $indexStyleDonor.Copy()
    $q3.Cells.Item($rowNum, 1).PasteSpecial($xlPasteFormats)

    # Columns B..G: text values, default (no) style.
    $q3.Cells.Item($rowNum, 2).Value = "'" + $rowData[0]   # 基金代码
    $q3.Cells.Item($rowNum, 2).Style = "Normal"
    $q3.Cells.Item($rowNum, 3).Value = "'" + $rowData[1]   # 基金名称
    $q3.Cells.Item($rowNum, 3).Style = "Normal"
    $q3.Cells.Item($rowNum, 4).Value = "'" + $rowData[2]   # 基金规模
    $q3.Cells.Item($rowNum, 4).Style = "Normal"
    $q3.Cells.Item($rowNum, 5).Value = "'" + $rowData[3]   # 股票总仓位
    $q3.Cells.Item($rowNum, 5).Style = "Normal"
    $q3.Cells.Item($rowNum, 6).Value = "'" + $rowData[4]   # 仓位占比
    $q3.Cells.Item($rowNum, 6).Style = "Normal"
    $q3.Cells.Item($rowNum, 7).Value = "'" + $rowData[5]   # 持有市值(亿元)
    $q3.Cells.Item($rowNum, 7).Style = "Normal"

    # Column H: numeric rank.
    $q3.Cells.Item($rowNum, 8).Value = $rowData[6]
}

# ---------------------------------------------------------------------------
# ---- 4. Update "总计": insert the 2022-Q3 row, shifting the rest down ----
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Snapshot rows 2..7 (the six existing quarters) before they're overwritten.
$existing = @()
for ($r = 2; $r -le 7; $r++) {
    $existing += ,@($zj.Cells.Item($r, 2).Value2, $zj.Cells.Item($r, 3).Value2, $zj.Cells.Item($r, 4).Value2)
}

$zjIndexDonor = $zj.Range("A7")

# Rewrite rows 3..8 from the snapshot (row N gets what used to be in row N-1).
for ($r = 8; $r -ge 3; $r--) {
    $src = $existing[$r - 3]

    $zj.Cells.Item($r, 1).Value = $r - 2
    $zjIndexDonor.Copy()
    $zj.Cells.Item($r, 1).PasteSpecial($xlPasteFormats)

    $zj.Cells.Item($r, 2).Value = "'" + $src[0]
    $zj.Cells.Item($r, 2).Style = "Normal"
    $zj.Cells.Item($r, 3).Value = $src[1]
    $zj.Cells.Item($r, 4).Value = $src[2]
}

# Row 2: the new 2022-Q3 summary values.
$zj.Cells.Item(2, 1).Value = 0
$zjIndexDonor.Copy()
$zj.Cells.Item(2, 1).PasteSpecial($xlPasteFormats)

$zj.Cells.Item(2, 2).Value = "'2022-Q3"
$zj.Cells.Item(2, 2).Style = "Normal"
$zj.Cells.Item(2, 3).Value = 16
$zj.Cells.Item(2, 4).Value = 5.19
